# Update the "想去人数" (interested-count) values in column F on both the
# "展览" and "全部类型" worksheets, which keep the same data in parallel.
#
# row -> new value for column F
$updates = @{
    4  = 47
    6  = 36
    7  = 120
    9  = 245
    10 = 14
    13 = 82
    14 = 353
    15 = 40
    17 = 391
    18 = 133
    19 = 60
    20 = 31
    22 = 962
    23 = 2743
    26 = 525
    27 = 966
    29 = 449
    31 = 384
    33 = 593
}

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
